$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.108.70'
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = '  -0.03%  '

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.563.73'
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = '  +0.70%  '

$ws.Range("E4").Value = '  -0.06%  '

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.56'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +2.66%  '

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.75'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +0.23%  '

$ws.Range("E8").Value = '  +1.51%  '

$ws.Range("E9").Value = '  +2.87%  '

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.64'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  +0.62%  '

$ws.Range("E11").Value = '  +0.00%  '

$ws.Range("E12").Value = '  +0.90%  '

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.37'
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  -0.49%  '

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.022.28'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  +0.61%  '

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.045.86'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  -0.09%  '

$ws.Range("E16").Value = '  +2.02%  '

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.524.44'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  -0.88%  '

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.34'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  -1.23%  '

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '343.85'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  +2.13%  '

$ws.Range("E20").Value = '  +2.82%  '

$ws.Range("E21").Value = '  +1.51%  '

$ws.Range("E22").Value = '  +0.02%  '

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.53'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  -3.69%  '

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.61'
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  +1.90%  '

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.694.62'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +0.59%  '

$ws.Range("E26").Value = '  +0.12%  '

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.62'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  -0.26%  '

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.12'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  +10.59%  '

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  +0.25%  '

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.49'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  -1.75%  '

$ws.Range("E31").Value = '  -0.57%  '

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.98'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  +6.64%  '

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0823'
$ws.Range("D33").Style = $origStyle

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '459.81'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  +11.59%  '

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '175.64'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  -0.42%  '

$ws.Range("E36").Value = '  +2.68%  '

$ws.Range("E37").Value = '  +2.05%  '

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.17'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  +0.72%  '

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.53'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +3.00%  '

$ws.Range("E40").Value = '  +0.05%  '

$ws.Range("E41").Value = '  -0.86%  '

$ws.Range("E42").Value = '  -0.05%  '

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '150.76'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  -1.73%  '

$ws.Range("E44").Value = '  +1.18%  '

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.84'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  -0.50%  '

$ws.Range("E46").Value = '  +4.42%  '

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.612'
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  +1.26%  '

$ws.Range("E48").Value = '  +1.17%  '

$ws.Range("E49").Value = '  +0.29%  '

$ws.Range("E50").Value = '  -2.46%  '

$ws.Range("E51").Value = '  +0.37%  '
